$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix name/text fields: the scraper had used commas as separators; replace with periods ---
$ws.Range('E24').Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range('E29').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('F29').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('E32').Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range('E45').Value = 'ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN'
$ws.Range('E94').Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range('E96').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('F96').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('E97').Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range('E114').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('F114').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('E126').Value = 'RICCOTTI. MARIANA EDITH'
$ws.Range('E130').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('F130').Value = 'IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA'
$ws.Range('E168').Value = 'MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO'
$ws.Range('E170').Value = 'DODERA. JORGE ABELARDO'
$ws.Range('E172').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'
$ws.Range('E178').Value = 'FERNANDEZ. MARIO HUGO'
$ws.Range('E189').Value = 'SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH'

# --- Fix "Importe" column (H2:H227): the scraper wrote Spanish-formatted numbers
# (thousands "." + decimal ",") as text; rewrite as plain dot-decimal text, keeping
# the cells as text (not auto-converted to actual numbers) via an explicit Text format. ---
$importeRange = $ws.Range("H2:H227")
$importeRange.NumberFormat = "@"
$ws.Range('H2').Value = '27780.80'
$ws.Range('H3').Value = '229000.00'
$ws.Range('H4').Value = '417401.60'
$ws.Range('H5').Value = '194150.86'
$ws.Range('H6').Value = '1684.90'
$ws.Range('H7').Value = '1939.83'
$ws.Range('H8').Value = '1316.39'
$ws.Range('H9').Value = '94080.00'
$ws.Range('H10').Value = '102205.80'
$ws.Range('H11').Value = '7340.00'
$ws.Range('H12').Value = '25169.99'
$ws.Range('H13').Value = '3535.00'
$ws.Range('H14').Value = '1209.60'
$ws.Range('H15').Value = '11616.15'
$ws.Range('H16').Value = '10702.51'
$ws.Range('H17').Value = '19620.00'
$ws.Range('H18').Value = '2337.50'
$ws.Range('H19').Value = '162.40'
$ws.Range('H20').Value = '12234.00'
$ws.Range('H21').Value = '170.00'
$ws.Range('H22').Value = '1706.10'
$ws.Range('H23').Value = '41280.00'
$ws.Range('H24').Value = '15.00'
$ws.Range('H25').Value = '41695.81'
$ws.Range('H26').Value = '205.00'
$ws.Range('H27').Value = '373.81'
$ws.Range('H28').Value = '1377.60'
$ws.Range('H29').Value = '47.85'
$ws.Range('H30').Value = '322.20'
$ws.Range('H31').Value = '762.00'
$ws.Range('H32').Value = '37.50'
$ws.Range('H33').Value = '5396.31'
$ws.Range('H34').Value = '20599.00'
$ws.Range('H35').Value = '96.00'
$ws.Range('H36').Value = '4968.00'
$ws.Range('H37').Value = '119.25'
$ws.Range('H38').Value = '19337.30'
$ws.Range('H39').Value = '279.00'
$ws.Range('H40').Value = '4078.90'
$ws.Range('H41').Value = '208.58'
$ws.Range('H42').Value = '935.00'
$ws.Range('H43').Value = '143421.60'
$ws.Range('H44').Value = '10.50'
$ws.Range('H45').Value = '1200.00'
$ws.Range('H46').Value = '36.00'
$ws.Range('H47').Value = '4309.20'
$ws.Range('H48').Value = '3791.78'
$ws.Range('H49').Value = '7726.10'
$ws.Range('H50').Value = '85700.00'
$ws.Range('H51').Value = '613.28'
$ws.Range('H52').Value = '130.00'
$ws.Range('H53').Value = '75.00'
$ws.Range('H54').Value = '27.82'
$ws.Range('H55').Value = '2546.00'
$ws.Range('H56').Value = '155.99'
$ws.Range('H57').Value = '1118.29'
$ws.Range('H58').Value = '170.72'
$ws.Range('H59').Value = '2700.00'
$ws.Range('H60').Value = '46607.86'
$ws.Range('H61').Value = '41211.55'
$ws.Range('H62').Value = '13180.00'
$ws.Range('H63').Value = '7014.61'
$ws.Range('H64').Value = '10072.20'
$ws.Range('H65').Value = '15900.00'
$ws.Range('H66').Value = '557.75'
$ws.Range('H67').Value = '81.20'
$ws.Range('H68').Value = '3941.62'
$ws.Range('H69').Value = '5515.00'
$ws.Range('H70').Value = '3699.50'
$ws.Range('H71').Value = '150.00'
$ws.Range('H72').Value = '504.18'
$ws.Range('H73').Value = '7728.60'
$ws.Range('H74').Value = '2207.00'
$ws.Range('H75').Value = '48.00'
$ws.Range('H76').Value = '3160.00'
$ws.Range('H77').Value = '13844.00'
$ws.Range('H78').Value = '69000.00'
$ws.Range('H79').Value = '7000.00'
$ws.Range('H80').Value = '611.50'
$ws.Range('H81').Value = '215.00'
$ws.Range('H82').Value = '2446.58'
$ws.Range('H83').Value = '5500.00'
$ws.Range('H84').Value = '960.00'
$ws.Range('H85').Value = '1064.00'
$ws.Range('H86').Value = '9983.00'
$ws.Range('H87').Value = '3737.60'
$ws.Range('H88').Value = '490.00'
$ws.Range('H89').Value = '455.00'
$ws.Range('H90').Value = '828.00'
$ws.Range('H91').Value = '500.00'
$ws.Range('H92').Value = '299.00'
$ws.Range('H93').Value = '1658.00'
$ws.Range('H94').Value = '92.00'
$ws.Range('H95').Value = '12675.30'
$ws.Range('H96').Value = '878.43'
$ws.Range('H97').Value = '3108.00'
$ws.Range('H98').Value = '20030.00'
$ws.Range('H99').Value = '285.00'
$ws.Range('H100').Value = '37650.00'
$ws.Range('H101').Value = '1450.00'
$ws.Range('H102').Value = '88.26'
$ws.Range('H103').Value = '31000.00'
$ws.Range('H104').Value = '194700.00'
$ws.Range('H105').Value = '8706.80'
$ws.Range('H106').Value = '6.03'
$ws.Range('H107').Value = '2490.00'
$ws.Range('H108').Value = '93.94'
$ws.Range('H109').Value = '4212.45'
$ws.Range('H110').Value = '1200.00'
$ws.Range('H111').Value = '3.68'
$ws.Range('H112').Value = '488.00'
$ws.Range('H113').Value = '133.75'
$ws.Range('H114').Value = '320.55'
$ws.Range('H115').Value = '347.50'
$ws.Range('H116').Value = '636.00'
$ws.Range('H117').Value = '295.20'
$ws.Range('H118').Value = '6750.01'
$ws.Range('H119').Value = '4118.12'
$ws.Range('H120').Value = '2328.04'
$ws.Range('H121').Value = '478.55'
$ws.Range('H122').Value = '114.74'
$ws.Range('H123').Value = '35851.00'
$ws.Range('H124').Value = '330.00'
$ws.Range('H125').Value = '1140.00'
$ws.Range('H126').Value = '3000.00'
$ws.Range('H127').Value = '640.00'
$ws.Range('H128').Value = '1200.00'
$ws.Range('H129').Value = '7.10'
$ws.Range('H130').Value = '1074.00'
$ws.Range('H131').Value = '7700.00'
$ws.Range('H132').Value = '25680.00'
$ws.Range('H133').Value = '3726.00'
$ws.Range('H134').Value = '1120.00'
$ws.Range('H135').Value = '1000.00'
$ws.Range('H136').Value = '604.64'
$ws.Range('H137').Value = '2631.87'
$ws.Range('H138').Value = '1152.50'
$ws.Range('H139').Value = '4212.20'
$ws.Range('H140').Value = '169.50'
$ws.Range('H141').Value = '1128.48'
$ws.Range('H142').Value = '24.04'
$ws.Range('H143').Value = '21055.00'
$ws.Range('H144').Value = '123.82'
$ws.Range('H145').Value = '130946.00'
$ws.Range('H146').Value = '7260.00'
$ws.Range('H147').Value = '1000.00'
$ws.Range('H148').Value = '1400.00'
$ws.Range('H149').Value = '1000.00'
$ws.Range('H150').Value = '12899.28'
$ws.Range('H151').Value = '384.00'
$ws.Range('H152').Value = '750.00'
$ws.Range('H153').Value = '2000.00'
$ws.Range('H154').Value = '4000.00'
$ws.Range('H155').Value = '20303.40'
$ws.Range('H156').Value = '1500.00'
$ws.Range('H157').Value = '950.00'
$ws.Range('H158').Value = '750.00'
$ws.Range('H159').Value = '6682.00'
$ws.Range('H160').Value = '2000.00'
$ws.Range('H161').Value = '1600.00'
$ws.Range('H162').Value = '200.00'
$ws.Range('H163').Value = '290.00'
$ws.Range('H164').Value = '7200.00'
$ws.Range('H165').Value = '2000.00'
$ws.Range('H166').Value = '400.00'
$ws.Range('H167').Value = '100.00'
$ws.Range('H168').Value = '80.00'
$ws.Range('H169').Value = '930.00'
$ws.Range('H170').Value = '1800.00'
$ws.Range('H171').Value = '4901.00'
$ws.Range('H172').Value = '2048.00'
$ws.Range('H173').Value = '2351.82'
$ws.Range('H174').Value = '2001.58'
$ws.Range('H175').Value = '3125.00'
$ws.Range('H176').Value = '1110.00'
$ws.Range('H177').Value = '160.00'
$ws.Range('H178').Value = '40.00'
$ws.Range('H179').Value = '718.62'
$ws.Range('H180').Value = '220.00'
$ws.Range('H181').Value = '876.84'
$ws.Range('H182').Value = '1386.50'
$ws.Range('H183').Value = '2043.88'
$ws.Range('H184').Value = '1458.54'
$ws.Range('H185').Value = '523.78'
$ws.Range('H186').Value = '1200.00'
$ws.Range('H187').Value = '135.40'
$ws.Range('H188').Value = '330.00'
$ws.Range('H189').Value = '21330.00'
$ws.Range('H190').Value = '6239.00'
$ws.Range('H191').Value = '1936.70'
$ws.Range('H192').Value = '990.00'
$ws.Range('H193').Value = '8355.68'
$ws.Range('H194').Value = '1770.00'
$ws.Range('H195').Value = '10150.00'
$ws.Range('H196').Value = '453.60'
$ws.Range('H197').Value = '12728.00'
$ws.Range('H198').Value = '5315.82'
$ws.Range('H199').Value = '2805.00'
$ws.Range('H200').Value = '22.35'
$ws.Range('H201').Value = '2052.60'
$ws.Range('H202').Value = '6016.69'
$ws.Range('H203').Value = '3405.02'
$ws.Range('H204').Value = '18000.00'
$ws.Range('H205').Value = '285.00'
$ws.Range('H206').Value = '2323.06'
$ws.Range('H207').Value = '23497.78'
$ws.Range('H208').Value = '510.00'
$ws.Range('H209').Value = '68200.00'
$ws.Range('H210').Value = '994.00'
$ws.Range('H211').Value = '136500.00'
$ws.Range('H212').Value = '62500.00'
$ws.Range('H213').Value = '74000.00'
$ws.Range('H214').Value = '20000.00'
$ws.Range('H215').Value = '223000.00'
$ws.Range('H216').Value = '32000.00'
$ws.Range('H217').Value = '111500.00'
$ws.Range('H218').Value = '223000.00'
$ws.Range('H219').Value = '1027139.98'
$ws.Range('H220').Value = '927193.79'
$ws.Range('H221').Value = '40000.00'
$ws.Range('H222').Value = '58000.00'
$ws.Range('H223').Value = '68000.00'
$ws.Range('H224').Value = '3550.00'
$ws.Range('H225').Value = '32000.00'
$ws.Range('H226').Value = '661.00'
$ws.Range('H227').Value = '1513.00'
